# Update the "Data" table with the latest Annual Population Survey (APS) release
# data, as described in the commit message:
#   "added in latest aps data added in version control added in note on aps accreditation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2, 3, 4 hold APS-sourced rows (Employment volumes / by occupation / by industry).
# Their "Latest period (release date)" (col C) and "Next period (release date)" (col D)
# values are rolled forward to the newest APS release.
$ws.Range("C2").Value = "Jul 2023 - Jun 2024 (15/10/24)"
$ws.Range("D2").Value = "Oct 2023 - Sep 2024 (14/01/25)"

$ws.Range("C3").Value = "Jul 2023 - Jun 2024 (15/10/24)"
$ws.Range("D3").Value = "Oct 2023 - Sep 2024 (14/01/25)"

$ws.Range("C4").Value = "Jul 2023 - Jun 2024 (15/10/24)"
$ws.Range("D4").Value = "Oct 2023 - Sep 2024 (14/01/25)"

# Update the view: move the selection to D5.
$ws.Range("D5").Select()
